$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# --- 1) Refresh the panel_query_time (F column) stamps on "data" ---
$ws.Range("F2").Value = "2021-10-05 14:22:19.899069"
$ws.Range("F3").Value = "2021-10-05 14:22:19.899077"
$ws.Range("F4").Value = "2021-10-05 14:22:19.899080"
$ws.Range("F5").Value = "2021-10-05 14:22:19.899082"
$ws.Range("F6").Value = "2021-10-05 14:22:19.899085"
$ws.Range("F7").Value = "2021-10-05 14:22:19.899087"
$ws.Range("F8").Value = "2021-10-05 14:22:19.899089"
$ws.Range("F9").Value = "2021-10-05 14:22:19.899091"
$ws.Range("F10").Value = "2021-10-05 14:22:19.899094"
$ws.Range("F11").Value = "2021-10-05 14:22:19.899096"
$ws.Range("F12").Value = "2021-10-05 14:22:19.899098"
$ws.Range("F13").Value = "2021-10-05 14:22:19.899100"
$ws.Range("F14").Value = "2021-10-05 14:22:19.899102"
$ws.Range("F15").Value = "2021-10-05 14:22:19.899104"
$ws.Range("F16").Value = "2021-10-05 14:22:19.899106"
$ws.Range("F17").Value = "2021-10-05 14:22:19.899109"
$ws.Range("F18").Value = "2021-10-05 14:22:19.899111"
$ws.Range("F19").Value = "2021-10-05 14:22:19.899113"
$ws.Range("F20").Value = "2021-10-05 14:22:19.899115"
$ws.Range("F21").Value = "2021-10-05 14:22:19.899118"
$ws.Range("F22").Value = "2021-10-05 14:22:19.899120"
$ws.Range("F23").Value = "2021-10-05 14:22:19.899122"
$ws.Range("F24").Value = "2021-10-05 14:22:19.899124"
$ws.Range("F25").Value = "2021-10-05 14:22:19.899126"
$ws.Range("F26").Value = "2021-10-05 14:22:19.899128"
$ws.Range("F27").Value = "2021-10-05 14:22:19.899131"
$ws.Range("F28").Value = "2021-10-05 14:22:19.899133"
$ws.Range("F29").Value = "2021-10-05 14:22:19.899135"
$ws.Range("F30").Value = "2021-10-05 14:22:19.899137"
$ws.Range("F31").Value = "2021-10-05 14:22:19.899139"
$ws.Range("F32").Value = "2021-10-05 14:22:19.899141"
$ws.Range("F33").Value = "2021-10-05 14:22:19.899143"
$ws.Range("F34").Value = "2021-10-05 14:22:19.899146"
$ws.Range("F35").Value = "2021-10-05 14:22:19.899148"
$ws.Range("F36").Value = "2021-10-05 14:22:19.899150"
$ws.Range("F37").Value = "2021-10-05 14:22:19.899153"
$ws.Range("F38").Value = "2021-10-05 14:22:19.899155"
$ws.Range("F39").Value = "2021-10-05 14:22:19.899157"
$ws.Range("F40").Value = "2021-10-05 14:22:19.899159"
$ws.Range("F41").Value = "2021-10-05 14:22:19.899161"
$ws.Range("F42").Value = "2021-10-05 14:22:19.899164"
$ws.Range("F43").Value = "2021-10-05 14:22:19.899166"
$ws.Range("F44").Value = "2021-10-05 14:22:19.899168"
$ws.Range("F45").Value = "2021-10-05 14:22:19.899170"
$ws.Range("F46").Value = "2021-10-05 14:22:19.899172"
$ws.Range("F47").Value = "2021-10-05 14:22:19.899174"
$ws.Range("F48").Value = "2021-10-05 14:22:19.899176"
$ws.Range("F49").Value = "2021-10-05 14:22:19.899179"
$ws.Range("F50").Value = "2021-10-05 14:22:19.899181"
$ws.Range("F51").Value = "2021-10-05 14:22:19.899183"
$ws.Range("F52").Value = "2021-10-05 14:22:19.899185"
$ws.Range("F53").Value = "2021-10-05 14:22:19.899187"
$ws.Range("F54").Value = "2021-10-05 14:22:19.899190"
$ws.Range("F55").Value = "2021-10-05 14:22:19.899192"
$ws.Range("F56").Value = "2021-10-05 14:22:19.899194"
$ws.Range("F57").Value = "2021-10-05 14:22:19.899196"
$ws.Range("F58").Value = "2021-10-05 14:22:19.899198"
$ws.Range("F59").Value = "2021-10-05 14:22:19.899201"
$ws.Range("F60").Value = "2021-10-05 14:22:19.899203"
$ws.Range("F61").Value = "2021-10-05 14:22:19.899205"
$ws.Range("F62").Value = "2021-10-05 14:22:19.899207"
$ws.Range("F63").Value = "2021-10-05 14:22:19.899209"
$ws.Range("F64").Value = "2021-10-05 14:22:19.899211"
$ws.Range("F65").Value = "2021-10-05 14:22:19.899213"

# --- 2) Add a new "metadata" worksheet positioned right after "data" ---
$meta = $wb.Worksheets.Add($null, $ws)
$meta.Name = "metadata"

# --- 3) Populate header row (styled like the "data" sheet headers) ---
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Copy the header style from the "data" sheet (B1, style index "1") onto the new headers
$ws.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 4) Populate the single metadata data row (row 2) ---
$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Primary ovarian insufficiency"
$meta.Range("C2").Value = 155
# D2 ("1.49") must stay literal text rather than be auto-coerced to a number;
# force text format, assign, then drop back to the default "Normal" style so
# no extra numFmt/style index lingers on the cell.
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "1.49"
$meta.Range("D2").Style = "Normal"
$meta.Range("E2").Value = "2021-09-28T12:11:55.816427Z"
$meta.Range("F2").Value = "2021-10-05 14:22:19.895969"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/155/?format=json"

# Copy the A2 style from the "data" sheet (style index "1") onto the new A2 cell
$ws.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 5) Leave selection on A1 of the metadata sheet, mirroring the source file ---
$meta.Range("A1").Select()

# Re-activate "data" so the workbook's active tab is unchanged by this edit
$ws.Activate()
